# Gnas-Tshr NATMI LR-pair sheet: refresh ligand/receptor expression and
# derived edge-weight metrics using newly recomputed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 7).Value = 108.9481836666667  # G2 Ligand average expression value
$ws.Cells.Item(2, 8).Value = 326.844551  # H2 Ligand total expression value
$ws.Cells.Item(2, 9).Value = 0.1523660837152667  # I2 Ligand derived specificity of average expression value
$ws.Cells.Item(2, 10).Value = 0.1650457680857909  # J2 Ligand derived specificity of total expression value
$ws.Cells.Item(2, 11).Value = 2  # K2 Receptor-expressing cells
$ws.Cells.Item(2, 12).Value = 0.6666666666666666  # L2 Receptor detection rate
$ws.Cells.Item(2, 13).Value = 0.1744923333333333  # M2 Receptor average expression value
$ws.Cells.Item(2, 14).Value = 0.523477  # N2 Receptor total expression value
$ws.Cells.Item(2, 15).Value = 0.07822917822503123  # O2 Receptor derived specificity of average expression value
$ws.Cells.Item(2, 16).Value = 0.08239975633156223  # P2 Receptor derived specificity of total expression value
$ws.Cells.Item(2, 17).Value = 19.01062278042522  # Q2 Edge average expression weight
$ws.Cells.Item(2, 18).Value = 171.095605023827  # R2 Edge total expression weight
$ws.Cells.Item(2, 19).Value = 0.01191947351841163  # S2 Edge average expression derived specificity
$ws.Cells.Item(2, 20).Value = 0.0135997310738247  # T2 Edge total expression derived specificity

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 7).Value = 108.9481836666667  # G3 Ligand average expression value
$ws.Cells.Item(3, 8).Value = 326.844551  # H3 Ligand total expression value
$ws.Cells.Item(3, 9).Value = 0.1523660837152667  # I3 Ligand derived specificity of average expression value
$ws.Cells.Item(3, 10).Value = 0.1650457680857909  # J3 Ligand derived specificity of total expression value
$ws.Cells.Item(3, 15).Value = 0.5425629886502931  # O3 Receptor derived specificity of average expression value
$ws.Cells.Item(3, 16).Value = 0.5714882742434749  # P3 Receptor derived specificity of total expression value
$ws.Cells.Item(3, 17).Value = 131.8492734537061  # Q3 Edge average expression weight
$ws.Cells.Item(3, 18).Value = 1186.643461083355  # R3 Edge total expression weight
$ws.Cells.Item(3, 19).Value = 0.08266819774949587  # S3 Edge average expression derived specificity
$ws.Cells.Item(3, 20).Value = 0.09432172117453745  # T3 Edge total expression derived specificity

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 7).Value = 108.9481836666667  # G4 Ligand average expression value
$ws.Cells.Item(4, 8).Value = 326.844551  # H4 Ligand total expression value
$ws.Cells.Item(4, 9).Value = 0.1523660837152667  # I4 Ligand derived specificity of average expression value
$ws.Cells.Item(4, 10).Value = 0.1650457680857909  # J4 Ligand derived specificity of total expression value
$ws.Cells.Item(4, 13).Value = 0.1427166666666667  # M4 Receptor average expression value
$ws.Cells.Item(4, 14).Value = 0.42815  # N4 Receptor total expression value
$ws.Cells.Item(4, 15).Value = 0.06398337015197826  # O4 Receptor derived specificity of average expression value
$ws.Cells.Item(4, 16).Value = 0.06739447133944447  # P4 Receptor derived specificity of total expression value
$ws.Cells.Item(4, 17).Value = 15.54872161229445  # Q4 Edge average expression weight
$ws.Cells.Item(4, 18).Value = 139.93849451065  # R4 Edge total expression weight
$ws.Cells.Item(4, 19).Value = 0.009748895532961217  # S4 Edge average expression derived specificity
$ws.Cells.Item(4, 20).Value = 0.01112317228695444  # T4 Edge total expression derived specificity

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 7).Value = 108.9481836666667  # G5 Ligand average expression value
$ws.Cells.Item(5, 8).Value = 326.844551  # H5 Ligand total expression value
$ws.Cells.Item(5, 9).Value = 0.1523660837152667  # I5 Ligand derived specificity of average expression value
$ws.Cells.Item(5, 10).Value = 0.1650457680857909  # J5 Ligand derived specificity of total expression value
$ws.Cells.Item(5, 13).Value = 0.3386875  # M5 Receptor average expression value
$ws.Cells.Item(5, 14).Value = 0.6773750000000001  # N5 Receptor total expression value
$ws.Cells.Item(5, 15).Value = 0.15184188493529  # O5 Receptor derived specificity of average expression value
$ws.Cells.Item(5, 16).Value = 0.1066246175956001  # P5 Receptor derived specificity of total expression value
$ws.Cells.Item(5, 17).Value = 36.89938795560418  # Q5 Edge average expression weight
$ws.Cells.Item(5, 18).Value = 221.396327733625  # R5 Edge total expression weight
$ws.Cells.Item(5, 19).Value = 0.02313555335153429  # S5 Edge average expression derived specificity
$ws.Cells.Item(5, 20).Value = 0.01759794190791956  # T5 Edge total expression derived specificity

# Row 6: ECs -> Resolving-Mac
$ws.Cells.Item(6, 7).Value = 108.9481836666667  # G6 Ligand average expression value
$ws.Cells.Item(6, 8).Value = 326.844551  # H6 Ligand total expression value
$ws.Cells.Item(6, 9).Value = 0.1523660837152667  # I6 Ligand derived specificity of average expression value
$ws.Cells.Item(6, 10).Value = 0.1650457680857909  # J6 Ligand derived specificity of total expression value
$ws.Cells.Item(6, 13).Value = 0.3644293333333333  # M6 Receptor average expression value
$ws.Cells.Item(6, 14).Value = 1.093288  # N6 Receptor total expression value
$ws.Cells.Item(6, 15).Value = 0.1633825780374074  # O6 Receptor derived specificity of average expression value
$ws.Cells.Item(6, 16).Value = 0.1720928804899184  # P6 Receptor derived specificity of total expression value
$ws.Cells.Item(6, 17).Value = 39.7039139415209  # Q6 Edge average expression weight
$ws.Cells.Item(6, 18).Value = 357.335225473688  # R6 Edge total expression weight
$ws.Cells.Item(6, 19).Value = 0.02489396356286372  # S6 Edge average expression derived specificity
$ws.Cells.Item(6, 20).Value = 0.02840320164255481  # T6 Edge total expression derived specificity

# Row 7: FAPs -> ECs
$ws.Cells.Item(7, 9).Value = 0.2954065074566193  # I7 Ligand derived specificity of average expression value
$ws.Cells.Item(7, 10).Value = 0.3199898083081954  # J7 Ligand derived specificity of total expression value
$ws.Cells.Item(7, 11).Value = 2  # K7 Receptor-expressing cells
$ws.Cells.Item(7, 12).Value = 0.6666666666666666  # L7 Receptor detection rate
$ws.Cells.Item(7, 13).Value = 0.1744923333333333  # M7 Receptor average expression value
$ws.Cells.Item(7, 14).Value = 0.523477  # N7 Receptor total expression value
$ws.Cells.Item(7, 15).Value = 0.07822917822503123  # O7 Receptor derived specificity of average expression value
$ws.Cells.Item(7, 16).Value = 0.08239975633156223  # P7 Receptor derived specificity of total expression value
$ws.Cells.Item(7, 17).Value = 36.85768868769555  # Q7 Edge average expression weight
$ws.Cells.Item(7, 18).Value = 331.7191981892599  # R7 Edge total expression weight
$ws.Cells.Item(7, 19).Value = 0.02310940832065789  # S7 Edge average expression derived specificity
$ws.Cells.Item(7, 20).Value = 0.02636708223317861  # T7 Edge total expression derived specificity

# Row 8: FAPs -> FAPs
$ws.Cells.Item(8, 9).Value = 0.2954065074566193  # I8 Ligand derived specificity of average expression value
$ws.Cells.Item(8, 10).Value = 0.3199898083081954  # J8 Ligand derived specificity of total expression value
$ws.Cells.Item(8, 15).Value = 0.5425629886502931  # O8 Receptor derived specificity of average expression value
$ws.Cells.Item(8, 16).Value = 0.5714882742434749  # P8 Receptor derived specificity of total expression value
$ws.Cells.Item(8, 19).Value = 0.1602766375524085  # S8 Edge average expression derived specificity
$ws.Cells.Item(8, 20).Value = 0.182870423325551  # T8 Edge total expression derived specificity

# Row 9: FAPs -> Inflammatory-Mac
$ws.Cells.Item(9, 9).Value = 0.2954065074566193  # I9 Ligand derived specificity of average expression value
$ws.Cells.Item(9, 10).Value = 0.3199898083081954  # J9 Ligand derived specificity of total expression value
$ws.Cells.Item(9, 13).Value = 0.1427166666666667  # M9 Receptor average expression value
$ws.Cells.Item(9, 14).Value = 0.42815  # N9 Receptor total expression value
$ws.Cells.Item(9, 15).Value = 0.06398337015197826  # O9 Receptor derived specificity of average expression value
$ws.Cells.Item(9, 16).Value = 0.06739447133944447  # P9 Receptor derived specificity of total expression value
$ws.Cells.Item(9, 17).Value = 30.14577414411111  # Q9 Edge average expression weight
$ws.Cells.Item(9, 18).Value = 271.311967297  # R9 Edge total expression weight
$ws.Cells.Item(9, 19).Value = 0.0189011039119  # S9 Edge average expression derived specificity
$ws.Cells.Item(9, 20).Value = 0.02156554396494101  # T9 Edge total expression derived specificity

# Row 10: FAPs -> MuSCs
$ws.Cells.Item(10, 9).Value = 0.2954065074566193  # I10 Ligand derived specificity of average expression value
$ws.Cells.Item(10, 10).Value = 0.3199898083081954  # J10 Ligand derived specificity of total expression value
$ws.Cells.Item(10, 13).Value = 0.3386875  # M10 Receptor average expression value
$ws.Cells.Item(10, 14).Value = 0.6773750000000001  # N10 Receptor total expression value
$ws.Cells.Item(10, 15).Value = 0.15184188493529  # O10 Receptor derived specificity of average expression value
$ws.Cells.Item(10, 16).Value = 0.1066246175956001  # P10 Receptor derived specificity of total expression value
$ws.Cells.Item(10, 17).Value = 71.54032615041666  # Q10 Edge average expression weight
$ws.Cells.Item(10, 18).Value = 429.2419569025  # R10 Edge total expression weight
$ws.Cells.Item(10, 19).Value = 0.04485508091436387  # S10 Edge average expression derived specificity
$ws.Cells.Item(10, 20).Value = 0.03411879094535073  # T10 Edge total expression derived specificity

# Row 11: FAPs -> Resolving-Mac
$ws.Cells.Item(11, 9).Value = 0.2954065074566193  # I11 Ligand derived specificity of average expression value
$ws.Cells.Item(11, 10).Value = 0.3199898083081954  # J11 Ligand derived specificity of total expression value
$ws.Cells.Item(11, 13).Value = 0.3644293333333333  # M11 Receptor average expression value
$ws.Cells.Item(11, 14).Value = 1.093288  # N11 Receptor total expression value
$ws.Cells.Item(11, 15).Value = 0.1633825780374074  # O11 Receptor derived specificity of average expression value
$ws.Cells.Item(11, 16).Value = 0.1720928804899184  # P11 Receptor derived specificity of total expression value
$ws.Cells.Item(11, 17).Value = 76.97772538238222  # Q11 Edge average expression weight
$ws.Cells.Item(11, 18).Value = 692.79952844144  # R11 Edge total expression weight
$ws.Cells.Item(11, 19).Value = 0.04826427675728908  # S11 Edge average expression derived specificity
$ws.Cells.Item(11, 20).Value = 0.05506796783917416  # T11 Edge total expression derived specificity

# Row 12: Inflammatory-Mac -> ECs
$ws.Cells.Item(12, 7).Value = 109.1710686666667  # G12 Ligand average expression value
$ws.Cells.Item(12, 8).Value = 327.513206  # H12 Ligand total expression value
$ws.Cells.Item(12, 9).Value = 0.1526777925792968  # I12 Ligand derived specificity of average expression value
$ws.Cells.Item(12, 10).Value = 0.1653834169091284  # J12 Ligand derived specificity of total expression value
$ws.Cells.Item(12, 11).Value = 2  # K12 Receptor-expressing cells
$ws.Cells.Item(12, 12).Value = 0.6666666666666666  # L12 Receptor detection rate
$ws.Cells.Item(12, 13).Value = 0.1744923333333333  # M12 Receptor average expression value
$ws.Cells.Item(12, 14).Value = 0.523477  # N12 Receptor total expression value
$ws.Cells.Item(12, 15).Value = 0.07822917822503123  # O12 Receptor derived specificity of average expression value
$ws.Cells.Item(12, 16).Value = 0.08239975633156223  # P12 Receptor derived specificity of total expression value
$ws.Cells.Item(12, 17).Value = 19.04951450414022  # Q12 Edge average expression weight
$ws.Cells.Item(12, 18).Value = 171.445630537262  # R12 Edge total expression weight
$ws.Cells.Item(12, 19).Value = 0.01194385824669016  # S12 Edge average expression derived specificity
$ws.Cells.Item(12, 20).Value = 0.01362755325459335  # T12 Edge total expression derived specificity

# Row 13: Inflammatory-Mac -> FAPs
$ws.Cells.Item(13, 7).Value = 109.1710686666667  # G13 Ligand average expression value
$ws.Cells.Item(13, 8).Value = 327.513206  # H13 Ligand total expression value
$ws.Cells.Item(13, 9).Value = 0.1526777925792968  # I13 Ligand derived specificity of average expression value
$ws.Cells.Item(13, 10).Value = 0.1653834169091284  # J13 Ligand derived specificity of total expression value
$ws.Cells.Item(13, 15).Value = 0.5425629886502931  # O13 Receptor derived specificity of average expression value
$ws.Cells.Item(13, 16).Value = 0.5714882742434749  # P13 Receptor derived specificity of total expression value
$ws.Cells.Item(13, 17).Value = 132.1190092521811  # Q13 Edge average expression weight
$ws.Cells.Item(13, 18).Value = 1189.07108326963  # R13 Edge total expression weight
$ws.Cells.Item(13, 19).Value = 0.08283731944235281  # S13 Edge average expression derived specificity
$ws.Cells.Item(13, 20).Value = 0.09451468351788692  # T13 Edge total expression derived specificity

# Row 14: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(14, 7).Value = 109.1710686666667  # G14 Ligand average expression value
$ws.Cells.Item(14, 8).Value = 327.513206  # H14 Ligand total expression value
$ws.Cells.Item(14, 9).Value = 0.1526777925792968  # I14 Ligand derived specificity of average expression value
$ws.Cells.Item(14, 10).Value = 0.1653834169091284  # J14 Ligand derived specificity of total expression value
$ws.Cells.Item(14, 13).Value = 0.1427166666666667  # M14 Receptor average expression value
$ws.Cells.Item(14, 14).Value = 0.42815  # N14 Receptor total expression value
$ws.Cells.Item(14, 15).Value = 0.06398337015197826  # O14 Receptor derived specificity of average expression value
$ws.Cells.Item(14, 16).Value = 0.06739447133944447  # P14 Receptor derived specificity of total expression value
$ws.Cells.Item(14, 17).Value = 15.58053101654445  # Q14 Edge average expression weight
$ws.Cells.Item(14, 18).Value = 140.2247791489  # R14 Edge total expression weight
$ws.Cells.Item(14, 19).Value = 0.009768839716588108  # S14 Edge average expression derived specificity
$ws.Cells.Item(14, 20).Value = 0.01114592795090165  # T14 Edge total expression derived specificity

# Row 15: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(15, 7).Value = 109.1710686666667  # G15 Ligand average expression value
$ws.Cells.Item(15, 8).Value = 327.513206  # H15 Ligand total expression value
$ws.Cells.Item(15, 9).Value = 0.1526777925792968  # I15 Ligand derived specificity of average expression value
$ws.Cells.Item(15, 10).Value = 0.1653834169091284  # J15 Ligand derived specificity of total expression value
$ws.Cells.Item(15, 13).Value = 0.3386875  # M15 Receptor average expression value
$ws.Cells.Item(15, 14).Value = 0.6773750000000001  # N15 Receptor total expression value
$ws.Cells.Item(15, 15).Value = 0.15184188493529  # O15 Receptor derived specificity of average expression value
$ws.Cells.Item(15, 16).Value = 0.1066246175956001  # P15 Receptor derived specificity of total expression value
$ws.Cells.Item(15, 17).Value = 36.97487631904167  # Q15 Edge average expression weight
$ws.Cells.Item(15, 18).Value = 221.84925791425  # R15 Edge total expression weight
$ws.Cells.Item(15, 19).Value = 0.02318288381299966  # S15 Edge average expression derived specificity
$ws.Cells.Item(15, 20).Value = 0.01763394358458953  # T15 Edge total expression derived specificity

# Row 16: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(16, 7).Value = 109.1710686666667  # G16 Ligand average expression value
$ws.Cells.Item(16, 8).Value = 327.513206  # H16 Ligand total expression value
$ws.Cells.Item(16, 9).Value = 0.1526777925792968  # I16 Ligand derived specificity of average expression value
$ws.Cells.Item(16, 10).Value = 0.1653834169091284  # J16 Ligand derived specificity of total expression value
$ws.Cells.Item(16, 13).Value = 0.3644293333333333  # M16 Receptor average expression value
$ws.Cells.Item(16, 14).Value = 1.093288  # N16 Receptor total expression value
$ws.Cells.Item(16, 15).Value = 0.1633825780374074  # O16 Receptor derived specificity of average expression value
$ws.Cells.Item(16, 16).Value = 0.1720928804899184  # P16 Receptor derived specificity of total expression value
$ws.Cells.Item(16, 17).Value = 39.78513977348089  # Q16 Edge average expression weight
$ws.Cells.Item(16, 18).Value = 358.0662579613281  # R16 Edge total expression weight
$ws.Cells.Item(16, 19).Value = 0.02494489136066607  # S16 Edge average expression derived specificity
$ws.Cells.Item(16, 20).Value = 0.02846130860115698  # T16 Edge total expression derived specificity

# Row 17: MuSCs -> ECs
$ws.Cells.Item(17, 7).Value = 164.799919  # G17 Ligand average expression value
$ws.Cells.Item(17, 8).Value = 329.599838  # H17 Ligand total expression value
$ws.Cells.Item(17, 9).Value = 0.2304757859153342  # I17 Ligand derived specificity of average expression value
$ws.Cells.Item(17, 10).Value = 0.166437097565877  # J17 Ligand derived specificity of total expression value
$ws.Cells.Item(17, 11).Value = 2  # K17 Receptor-expressing cells
$ws.Cells.Item(17, 12).Value = 0.6666666666666666  # L17 Receptor detection rate
$ws.Cells.Item(17, 13).Value = 0.1744923333333333  # M17 Receptor average expression value
$ws.Cells.Item(17, 14).Value = 0.523477  # N17 Receptor total expression value
$ws.Cells.Item(17, 15).Value = 0.07822917822503123  # O17 Receptor derived specificity of average expression value
$ws.Cells.Item(17, 16).Value = 0.08239975633156223  # P17 Receptor derived specificity of total expression value
$ws.Cells.Item(17, 17).Value = 28.75632239945433  # Q17 Edge average expression weight
$ws.Cells.Item(17, 18).Value = 172.537934396726  # R17 Edge total expression weight
$ws.Cells.Item(17, 19).Value = 0.01802993133292482  # S17 Edge average expression derived specificity
$ws.Cells.Item(17, 20).Value = 0.01371437628396071  # T17 Edge total expression derived specificity

# Row 18: MuSCs -> FAPs
$ws.Cells.Item(18, 7).Value = 164.799919  # G18 Ligand average expression value
$ws.Cells.Item(18, 8).Value = 329.599838  # H18 Ligand total expression value
$ws.Cells.Item(18, 9).Value = 0.2304757859153342  # I18 Ligand derived specificity of average expression value
$ws.Cells.Item(18, 10).Value = 0.166437097565877  # J18 Ligand derived specificity of total expression value
$ws.Cells.Item(18, 15).Value = 0.5425629886502931  # O18 Receptor derived specificity of average expression value
$ws.Cells.Item(18, 16).Value = 0.5714882742434749  # P18 Receptor derived specificity of total expression value
$ws.Cells.Item(18, 17).Value = 199.4411366403316  # Q18 Edge average expression weight
$ws.Cells.Item(18, 18).Value = 1196.64681984199  # R18 Edge total expression weight
$ws.Cells.Item(18, 19).Value = 0.1250476312177488  # S18 Edge average expression derived specificity
$ws.Cells.Item(18, 20).Value = 0.09511684965801591  # T18 Edge total expression derived specificity

# Row 19: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(19, 7).Value = 164.799919  # G19 Ligand average expression value
$ws.Cells.Item(19, 8).Value = 329.599838  # H19 Ligand total expression value
$ws.Cells.Item(19, 9).Value = 0.2304757859153342  # I19 Ligand derived specificity of average expression value
$ws.Cells.Item(19, 10).Value = 0.166437097565877  # J19 Ligand derived specificity of total expression value
$ws.Cells.Item(19, 13).Value = 0.1427166666666667  # M19 Receptor average expression value
$ws.Cells.Item(19, 14).Value = 0.42815  # N19 Receptor total expression value
$ws.Cells.Item(19, 15).Value = 0.06398337015197826  # O19 Receptor derived specificity of average expression value
$ws.Cells.Item(19, 16).Value = 0.06739447133944447  # P19 Receptor derived specificity of total expression value
$ws.Cells.Item(19, 17).Value = 23.51969510661667  # Q19 Edge average expression weight
$ws.Cells.Item(19, 18).Value = 141.1181706397  # R19 Edge total expression weight
$ws.Cells.Item(19, 19).Value = 0.01474661752128892  # S19 Edge average expression derived specificity
$ws.Cells.Item(19, 20).Value = 0.01121694020172382  # T19 Edge total expression derived specificity

# Row 20: MuSCs -> MuSCs
$ws.Cells.Item(20, 7).Value = 164.799919  # G20 Ligand average expression value
$ws.Cells.Item(20, 8).Value = 329.599838  # H20 Ligand total expression value
$ws.Cells.Item(20, 9).Value = 0.2304757859153342  # I20 Ligand derived specificity of average expression value
$ws.Cells.Item(20, 10).Value = 0.166437097565877  # J20 Ligand derived specificity of total expression value
$ws.Cells.Item(20, 13).Value = 0.3386875  # M20 Receptor average expression value
$ws.Cells.Item(20, 14).Value = 0.6773750000000001  # N20 Receptor total expression value
$ws.Cells.Item(20, 15).Value = 0.15184188493529  # O20 Receptor derived specificity of average expression value
$ws.Cells.Item(20, 16).Value = 0.1066246175956001  # P20 Receptor derived specificity of total expression value
$ws.Cells.Item(20, 17).Value = 55.8156725663125  # Q20 Edge average expression weight
$ws.Cells.Item(20, 18).Value = 223.26269026525  # R20 Edge total expression weight
$ws.Cells.Item(20, 19).Value = 0.0349958777653267  # S20 Edge average expression derived specificity
$ws.Cells.Item(20, 20).Value = 0.01774629188168323  # T20 Edge total expression derived specificity

# Row 21: MuSCs -> Resolving-Mac
$ws.Cells.Item(21, 7).Value = 164.799919  # G21 Ligand average expression value
$ws.Cells.Item(21, 8).Value = 329.599838  # H21 Ligand total expression value
$ws.Cells.Item(21, 9).Value = 0.2304757859153342  # I21 Ligand derived specificity of average expression value
$ws.Cells.Item(21, 10).Value = 0.166437097565877  # J21 Ligand derived specificity of total expression value
$ws.Cells.Item(21, 13).Value = 0.3644293333333333  # M21 Receptor average expression value
$ws.Cells.Item(21, 14).Value = 1.093288  # N21 Receptor total expression value
$ws.Cells.Item(21, 15).Value = 0.1633825780374074  # O21 Receptor derived specificity of average expression value
$ws.Cells.Item(21, 16).Value = 0.1720928804899184  # P21 Receptor derived specificity of total expression value
$ws.Cells.Item(21, 17).Value = 60.05792461455733  # Q21 Edge average expression weight
$ws.Cells.Item(21, 18).Value = 360.347547687344  # R21 Edge total expression weight
$ws.Cells.Item(21, 19).Value = 0.0376557280780449  # S21 Edge average expression derived specificity
$ws.Cells.Item(21, 20).Value = 0.02864263954049336  # T21 Edge total expression derived specificity

# Row 22: Resolving-Mac -> ECs
$ws.Cells.Item(22, 7).Value = 120.894928  # G22 Ligand average expression value
$ws.Cells.Item(22, 8).Value = 362.684784  # H22 Ligand total expression value
$ws.Cells.Item(22, 9).Value = 0.169073830333483  # I22 Ligand derived specificity of average expression value
$ws.Cells.Item(22, 10).Value = 0.1831439091310082  # J22 Ligand derived specificity of total expression value
$ws.Cells.Item(22, 11).Value = 2  # K22 Receptor-expressing cells
$ws.Cells.Item(22, 12).Value = 0.6666666666666666  # L22 Receptor detection rate
$ws.Cells.Item(22, 13).Value = 0.1744923333333333  # M22 Receptor average expression value
$ws.Cells.Item(22, 14).Value = 0.523477  # N22 Receptor total expression value
$ws.Cells.Item(22, 15).Value = 0.07822917822503123  # O22 Receptor derived specificity of average expression value
$ws.Cells.Item(22, 16).Value = 0.08239975633156223  # P22 Receptor derived specificity of total expression value
$ws.Cells.Item(22, 17).Value = 21.09523807488533  # Q22 Edge average expression weight
$ws.Cells.Item(22, 18).Value = 189.857142673968  # R22 Edge total expression weight
$ws.Cells.Item(22, 19).Value = 0.01322650680634673  # S22 Edge average expression derived specificity
$ws.Cells.Item(22, 20).Value = 0.01509101348600485  # T22 Edge total expression derived specificity

# Row 23: Resolving-Mac -> FAPs
$ws.Cells.Item(23, 7).Value = 120.894928  # G23 Ligand average expression value
$ws.Cells.Item(23, 8).Value = 362.684784  # H23 Ligand total expression value
$ws.Cells.Item(23, 9).Value = 0.169073830333483  # I23 Ligand derived specificity of average expression value
$ws.Cells.Item(23, 10).Value = 0.1831439091310082  # J23 Ligand derived specificity of total expression value
$ws.Cells.Item(23, 15).Value = 0.5425629886502931  # O23 Receptor derived specificity of average expression value
$ws.Cells.Item(23, 16).Value = 0.5714882742434749  # P23 Receptor derived specificity of total expression value
$ws.Cells.Item(23, 17).Value = 146.3072433571467  # Q23 Edge average expression weight
$ws.Cells.Item(23, 18).Value = 1316.76519021432  # R23 Edge total expression weight
$ws.Cells.Item(23, 19).Value = 0.09173320268828711  # S23 Edge average expression derived specificity
$ws.Cells.Item(23, 20).Value = 0.1046645965674837  # T23 Edge total expression derived specificity

# Row 24: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(24, 7).Value = 120.894928  # G24 Ligand average expression value
$ws.Cells.Item(24, 8).Value = 362.684784  # H24 Ligand total expression value
$ws.Cells.Item(24, 9).Value = 0.169073830333483  # I24 Ligand derived specificity of average expression value
$ws.Cells.Item(24, 10).Value = 0.1831439091310082  # J24 Ligand derived specificity of total expression value
$ws.Cells.Item(24, 13).Value = 0.1427166666666667  # M24 Receptor average expression value
$ws.Cells.Item(24, 14).Value = 0.42815  # N24 Receptor total expression value
$ws.Cells.Item(24, 15).Value = 0.06398337015197826  # O24 Receptor derived specificity of average expression value
$ws.Cells.Item(24, 16).Value = 0.06739447133944447  # P24 Receptor derived specificity of total expression value
$ws.Cells.Item(24, 17).Value = 17.25372114106667  # Q24 Edge average expression weight
$ws.Cells.Item(24, 18).Value = 155.2834902696  # R24 Edge total expression weight
$ws.Cells.Item(24, 19).Value = 0.01081791346924001  # S24 Edge average expression derived specificity
$ws.Cells.Item(24, 20).Value = 0.01234288693492356  # T24 Edge total expression derived specificity

# Row 25: Resolving-Mac -> MuSCs
$ws.Cells.Item(25, 7).Value = 120.894928  # G25 Ligand average expression value
$ws.Cells.Item(25, 8).Value = 362.684784  # H25 Ligand total expression value
$ws.Cells.Item(25, 9).Value = 0.169073830333483  # I25 Ligand derived specificity of average expression value
$ws.Cells.Item(25, 10).Value = 0.1831439091310082  # J25 Ligand derived specificity of total expression value
$ws.Cells.Item(25, 13).Value = 0.3386875  # M25 Receptor average expression value
$ws.Cells.Item(25, 14).Value = 0.6773750000000001  # N25 Receptor total expression value
$ws.Cells.Item(25, 15).Value = 0.15184188493529  # O25 Receptor derived specificity of average expression value
$ws.Cells.Item(25, 16).Value = 0.1066246175956001  # P25 Receptor derived specificity of total expression value
$ws.Cells.Item(25, 17).Value = 40.94560092700001  # Q25 Edge average expression weight
$ws.Cells.Item(25, 18).Value = 245.673605562  # R25 Edge total expression weight
$ws.Cells.Item(25, 19).Value = 0.02567248909106547  # S25 Edge average expression derived specificity
$ws.Cells.Item(25, 20).Value = 0.01952764927605709  # T25 Edge total expression derived specificity

# Row 26: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(26, 7).Value = 120.894928  # G26 Ligand average expression value
$ws.Cells.Item(26, 8).Value = 362.684784  # H26 Ligand total expression value
$ws.Cells.Item(26, 9).Value = 0.169073830333483  # I26 Ligand derived specificity of average expression value
$ws.Cells.Item(26, 10).Value = 0.1831439091310082  # J26 Ligand derived specificity of total expression value
$ws.Cells.Item(26, 11).Value = 2  # K26 Receptor-expressing cells
$ws.Cells.Item(26, 12).Value = 0.6666666666666666  # L26 Receptor detection rate
$ws.Cells.Item(26, 13).Value = 0.3644293333333333  # M26 Receptor average expression value
$ws.Cells.Item(26, 14).Value = 1.093288  # N26 Receptor total expression value
$ws.Cells.Item(26, 15).Value = 0.1633825780374074  # O26 Receptor derived specificity of average expression value
$ws.Cells.Item(26, 16).Value = 0.1720928804899184  # P26 Receptor derived specificity of total expression value
$ws.Cells.Item(26, 17).Value = 44.05765801442134  # Q26 Edge average expression weight
$ws.Cells.Item(26, 18).Value = 396.518922129792  # R26 Edge total expression weight
$ws.Cells.Item(26, 19).Value = 0.02762371827854367  # S26 Edge average expression derived specificity
$ws.Cells.Item(26, 20).Value = 0.03151776286653907  # T26 Edge total expression derived specificity

